# Apply updated crypto price/volume snapshot values (GitHub Actions scrape)
# matching commit "Updated cryptos list on Sat Mar  4 14:36:17 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.390.19'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '1.571.27'
$ws.Range('E3').Value = '  -0.25%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = "'1.002"
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').Value = "'291.34"
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('D7').Value = "'0.3764"
$ws.Range('E7').Value = '  +2.29%  '
$ws.Range('D8').Value = "'49.75"
$ws.Range('E8').Value = '  +0.74%  '
$ws.Range('D9').Value = "'0.3417"
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('D10').Value = "'0.07627"
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').Value = "'1.144"
$ws.Range('E11').Value = '  -2.29%  '
$ws.Range('D12').Value = "'1.003"
$ws.Range('E12').Value = '  +0.21%  '
$ws.Range('D13').Value = "'21.16"
$ws.Range('E13').Value = '  -1.02%  '
$ws.Range('D14').Value = "'6.003"
$ws.Range('E14').Value = '  -1.06%  '
$ws.Range('D15').Value = "'6.934"
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('D16').Value = '1.571.43'
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('D18').Value = "'90.30"
$ws.Range('E18').Value = '  +1.16%  '
$ws.Range('D19').Value = "'0.06742"
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('D21').Value = "'16.72"
$ws.Range('E21').Value = '  +0.71%  '
$ws.Range('D22').Value = "'6.208"
$ws.Range('E22').Value = '  -0.80%  '
$ws.Range('D23').Value = "'12.01"
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = '22.389.07'
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').Value = "'2.394"
$ws.Range('E25').Value = '  +0.45%  '
$ws.Range('D26').Value = "'2.655"
$ws.Range('E26').Value = '  -11.21%  '
$ws.Range('E27').Value = '  +0.53%  '
$ws.Range('D28').Value = "'147.02"
$ws.Range('E28').Value = '  +0.83%  '
$ws.Range('D29').Value = "'5.040"
$ws.Range('E29').Value = '  +1.20%  '
$ws.Range('D30').Value = "'126.75"
$ws.Range('E30').Value = '  +0.61%  '
$ws.Range('D31').Value = '1.748.70'
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('B32').Value = 'WEMIXTOKEN'
$ws.Range('C32').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D32').Value = "'2.012"
$ws.Range('E32').Value = '  +0.58%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'6.121"
$ws.Range('E33').Value = '  -2.97%  '
$ws.Range('D34').Value = "'0.9851"
$ws.Range('E34').Value = '  -6.00%  '
$ws.Range('D35').Value = "'10.14"
$ws.Range('E35').Value = '  -2.13%  '
$ws.Range('D36').Value = "'0.08483"
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('D37').Value = "'0.02536"
$ws.Range('E37').Value = '  -0.31%  '
$ws.Range('D38').Value = "'1.381"
$ws.Range('E38').Value = '  +10.63%  '
$ws.Range('D39').Value = "'0.2311"
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('D40').Value = "'0.06504"
$ws.Range('E40').Value = '  -1.39%  '
$ws.Range('D41').Value = "'5.428"
$ws.Range('E41').Value = '  -2.47%  '
$ws.Range('D42').Value = "'0.6353"
$ws.Range('E42').Value = '  -0.62%  '
$ws.Range('D43').Value = "'11.38"
$ws.Range('E43').Value = '  -3.78%  '
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('D45').Value = "'14.06"
$ws.Range('E45').Value = '  -2.81%  '
$ws.Range('D46').Value = "'3.792"
$ws.Range('D47').Value = "'0.5962"
$ws.Range('E47').Value = '  -0.79%  '
$ws.Range('D48').Value = "'2.093"
$ws.Range('E48').Value = '  -1.78%  '
$ws.Range('D49').Value = "'1.281"
$ws.Range('E49').Value = '  +1.63%  '
$ws.Range('D50').Value = "'124.57"
$ws.Range('E50').Value = '  +0.89%  '
